$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Update header metadata (row 2: last update time, row 3: total rows) on sheet1
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 16:36:34'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 362'

# Update "Ultima actualizacion" timestamp on sheet2 and sheet3 (row 2)
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 16:36:34'
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 16:36:34'

# Update dimension-affecting data rows (re-sorted schedule rows + 4 newly scraped rows)
$ws1.Cells.Item(50, 1).Value = '06:45:50'
$ws1.Cells.Item(50, 2).Value = '08:01'
$ws1.Cells.Item(50, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(50, 4).Value = 76
$ws1.Cells.Item(50, 5).Value = 'LP1912'
$ws1.Cells.Item(51, 1).Value = '07:12:53'
$ws1.Cells.Item(51, 2).Value = '08:01'
$ws1.Cells.Item(51, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(51, 4).Value = 49
$ws1.Cells.Item(51, 5).Value = 'LP1912'
$ws1.Cells.Item(84, 1).Value = '08:39:08'
$ws1.Cells.Item(84, 2).Value = '09:04'
$ws1.Cells.Item(84, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(84, 4).Value = 25
$ws1.Cells.Item(84, 5).Value = 'LP1912'
$ws1.Cells.Item(85, 1).Value = '07:36:59'
$ws1.Cells.Item(85, 2).Value = '09:04'
$ws1.Cells.Item(85, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(85, 4).Value = 88
$ws1.Cells.Item(85, 5).Value = 'LP1912'
$ws1.Cells.Item(107, 1).Value = '08:11:27'
$ws1.Cells.Item(107, 2).Value = '10:03'
$ws1.Cells.Item(107, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(107, 4).Value = 112
$ws1.Cells.Item(107, 5).Value = 'LP1912'
$ws1.Cells.Item(108, 1).Value = '09:21:49'
$ws1.Cells.Item(108, 2).Value = '10:03'
$ws1.Cells.Item(108, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(108, 4).Value = 42
$ws1.Cells.Item(108, 5).Value = 'LP1912'
$ws1.Cells.Item(118, 1).Value = '08:29:19'
$ws1.Cells.Item(118, 2).Value = '10:15'
$ws1.Cells.Item(118, 3).Value = '17_ROMERO'
$ws1.Cells.Item(118, 4).Value = 106
$ws1.Cells.Item(118, 5).Value = 'LP1912'
$ws1.Cells.Item(119, 1).Value = '08:39:08'
$ws1.Cells.Item(119, 2).Value = '10:15'
$ws1.Cells.Item(119, 3).Value = '10_OLMOS'
$ws1.Cells.Item(119, 4).Value = 96
$ws1.Cells.Item(119, 5).Value = 'LP1912'
$ws1.Cells.Item(130, 1).Value = '10:36:18'
$ws1.Cells.Item(130, 2).Value = '10:37'
$ws1.Cells.Item(130, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(130, 4).Value = 1
$ws1.Cells.Item(130, 5).Value = 'LP1912'
$ws1.Cells.Item(131, 1).Value = '08:39:08'
$ws1.Cells.Item(131, 2).Value = '10:37'
$ws1.Cells.Item(131, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(131, 4).Value = 118
$ws1.Cells.Item(131, 5).Value = 'LP1912'
$ws1.Cells.Item(140, 1).Value = '10:55:25'
$ws1.Cells.Item(140, 2).Value = '10:56'
$ws1.Cells.Item(140, 3).Value = '10_OLMOS'
$ws1.Cells.Item(140, 4).Value = 1
$ws1.Cells.Item(140, 5).Value = 'LP1912'
$ws1.Cells.Item(142, 1).Value = '10:55:25'
$ws1.Cells.Item(142, 2).Value = '10:56'
$ws1.Cells.Item(142, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(142, 4).Value = 1
$ws1.Cells.Item(142, 5).Value = 'LP1912'
$ws1.Cells.Item(164, 1).Value = '10:04:17'
$ws1.Cells.Item(164, 2).Value = '11:34'
$ws1.Cells.Item(164, 3).Value = '10_OLMOS'
$ws1.Cells.Item(164, 4).Value = 90
$ws1.Cells.Item(164, 5).Value = 'LP1912'
$ws1.Cells.Item(165, 1).Value = '11:34:25'
$ws1.Cells.Item(165, 2).Value = '11:34'
$ws1.Cells.Item(165, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(165, 4).Value = 0
$ws1.Cells.Item(165, 5).Value = 'LP1912'
$ws1.Cells.Item(166, 1).Value = '11:34:25'
$ws1.Cells.Item(166, 2).Value = '11:34'
$ws1.Cells.Item(166, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(166, 4).Value = 0
$ws1.Cells.Item(166, 5).Value = 'LP1912'
$ws1.Cells.Item(186, 1).Value = '12:11:45'
$ws1.Cells.Item(186, 2).Value = '12:17'
$ws1.Cells.Item(186, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(186, 4).Value = 6
$ws1.Cells.Item(186, 5).Value = 'LP1912'
$ws1.Cells.Item(187, 1).Value = '11:11:31'
$ws1.Cells.Item(187, 2).Value = '12:17'
$ws1.Cells.Item(187, 3).Value = '15_ABASTO'
$ws1.Cells.Item(187, 4).Value = 66
$ws1.Cells.Item(187, 5).Value = 'LP1912'
$ws1.Cells.Item(188, 1).Value = '11:53:59'
$ws1.Cells.Item(188, 2).Value = '12:17'
$ws1.Cells.Item(188, 3).Value = '10_OLMOS'
$ws1.Cells.Item(188, 4).Value = 24
$ws1.Cells.Item(188, 5).Value = 'LP1912'
$ws1.Cells.Item(220, 1).Value = '12:11:45'
$ws1.Cells.Item(220, 2).Value = '13:06'
$ws1.Cells.Item(220, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(220, 4).Value = 55
$ws1.Cells.Item(220, 5).Value = 'LP1912'
$ws1.Cells.Item(221, 1).Value = '11:11:31'
$ws1.Cells.Item(221, 2).Value = '13:06'
$ws1.Cells.Item(221, 3).Value = '14_ABASTO'
$ws1.Cells.Item(221, 4).Value = 115
$ws1.Cells.Item(221, 5).Value = 'LP1912'
$ws1.Cells.Item(232, 1).Value = '11:53:59'
$ws1.Cells.Item(232, 2).Value = '13:20'
$ws1.Cells.Item(232, 3).Value = '17_ROMERO'
$ws1.Cells.Item(232, 4).Value = 87
$ws1.Cells.Item(232, 5).Value = 'LP1912'
$ws1.Cells.Item(233, 1).Value = '11:53:59'
$ws1.Cells.Item(233, 2).Value = '13:20'
$ws1.Cells.Item(233, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(233, 4).Value = 87
$ws1.Cells.Item(233, 5).Value = 'LP1912'
$ws1.Cells.Item(273, 1).Value = '14:31:57'
$ws1.Cells.Item(273, 2).Value = '14:31'
$ws1.Cells.Item(273, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(273, 4).Value = 0
$ws1.Cells.Item(273, 5).Value = 'LP1912'
$ws1.Cells.Item(274, 1).Value = '13:12:59'
$ws1.Cells.Item(274, 2).Value = '14:31'
$ws1.Cells.Item(274, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(274, 4).Value = 79
$ws1.Cells.Item(274, 5).Value = 'LP1912'
$ws1.Cells.Item(275, 1).Value = '14:31:57'
$ws1.Cells.Item(275, 2).Value = '14:31'
$ws1.Cells.Item(275, 3).Value = '15_ABASTO'
$ws1.Cells.Item(275, 4).Value = 0
$ws1.Cells.Item(275, 5).Value = 'LP1912'
$ws1.Cells.Item(337, 1).Value = '16:36:34'
$ws1.Cells.Item(337, 2).Value = '16:44'
$ws1.Cells.Item(337, 3).Value = '14_ABASTO'
$ws1.Cells.Item(337, 4).Value = 8
$ws1.Cells.Item(337, 5).Value = 'LP1912'
$ws1.Cells.Item(338, 1).Value = '15:45:31'
$ws1.Cells.Item(338, 2).Value = '16:48'
$ws1.Cells.Item(338, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(338, 4).Value = 63
$ws1.Cells.Item(338, 5).Value = 'LP1912'
$ws1.Cells.Item(339, 1).Value = '15:57:19'
$ws1.Cells.Item(339, 2).Value = '16:51'
$ws1.Cells.Item(339, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(339, 4).Value = 54
$ws1.Cells.Item(339, 5).Value = 'LP1912'
$ws1.Cells.Item(340, 1).Value = '15:45:31'
$ws1.Cells.Item(340, 2).Value = '16:53'
$ws1.Cells.Item(340, 3).Value = '10_OLMOS'
$ws1.Cells.Item(340, 4).Value = 68
$ws1.Cells.Item(340, 5).Value = 'LP1912'
$ws1.Cells.Item(341, 1).Value = '15:57:19'
$ws1.Cells.Item(341, 2).Value = '16:54'
$ws1.Cells.Item(341, 3).Value = '10_OLMOS'
$ws1.Cells.Item(341, 4).Value = 57
$ws1.Cells.Item(341, 5).Value = 'LP1912'
$ws1.Cells.Item(342, 1).Value = '16:27:37'
$ws1.Cells.Item(342, 2).Value = '16:56'
$ws1.Cells.Item(342, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(342, 4).Value = 29
$ws1.Cells.Item(342, 5).Value = 'LP1912'
$ws1.Cells.Item(343, 1).Value = '15:17:21'
$ws1.Cells.Item(343, 2).Value = '16:57'
$ws1.Cells.Item(343, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(343, 4).Value = 100
$ws1.Cells.Item(343, 5).Value = 'LP1912'
$ws1.Cells.Item(344, 1).Value = '15:57:19'
$ws1.Cells.Item(344, 2).Value = '17:01'
$ws1.Cells.Item(344, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(344, 4).Value = 64
$ws1.Cells.Item(344, 5).Value = 'LP1912'
$ws1.Cells.Item(345, 1).Value = '15:17:21'
$ws1.Cells.Item(345, 2).Value = '17:04'
$ws1.Cells.Item(345, 3).Value = '14_ABASTO'
$ws1.Cells.Item(345, 4).Value = 107
$ws1.Cells.Item(345, 5).Value = 'LP1912'
$ws1.Cells.Item(346, 1).Value = '15:45:31'
$ws1.Cells.Item(346, 2).Value = '17:07'
$ws1.Cells.Item(346, 3).Value = '15_ABASTO'
$ws1.Cells.Item(346, 4).Value = 82
$ws1.Cells.Item(346, 5).Value = 'LP1912'
$ws1.Cells.Item(347, 1).Value = '16:27:37'
$ws1.Cells.Item(347, 2).Value = '17:14'
$ws1.Cells.Item(347, 3).Value = '10_OLMOS'
$ws1.Cells.Item(347, 4).Value = 47
$ws1.Cells.Item(347, 5).Value = 'LP1912'
$ws1.Cells.Item(348, 1).Value = '15:45:31'
$ws1.Cells.Item(348, 2).Value = '17:17'
$ws1.Cells.Item(348, 3).Value = '17_ROMERO'
$ws1.Cells.Item(348, 4).Value = 92
$ws1.Cells.Item(348, 5).Value = 'LP1912'
$ws1.Cells.Item(349, 1).Value = '15:45:31'
$ws1.Cells.Item(349, 2).Value = '17:24'
$ws1.Cells.Item(349, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(349, 4).Value = 99
$ws1.Cells.Item(349, 5).Value = 'LP1912'
$ws1.Cells.Item(350, 1).Value = '15:57:19'
$ws1.Cells.Item(350, 2).Value = '17:27'
$ws1.Cells.Item(350, 3).Value = '15_ABASTO'
$ws1.Cells.Item(350, 4).Value = 90
$ws1.Cells.Item(350, 5).Value = 'LP1912'
$ws1.Cells.Item(351, 1).Value = '16:36:34'
$ws1.Cells.Item(351, 2).Value = '17:30'
$ws1.Cells.Item(351, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(351, 4).Value = 54
$ws1.Cells.Item(351, 5).Value = 'LP1912'
$ws1.Cells.Item(352, 1).Value = '16:13:19'
$ws1.Cells.Item(352, 2).Value = '17:31'
$ws1.Cells.Item(352, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(352, 4).Value = 78
$ws1.Cells.Item(352, 5).Value = 'LP1912'
$ws1.Cells.Item(353, 1).Value = '15:45:31'
$ws1.Cells.Item(353, 2).Value = '17:34'
$ws1.Cells.Item(353, 3).Value = '10_OLMOS'
$ws1.Cells.Item(353, 4).Value = 109
$ws1.Cells.Item(353, 5).Value = 'LP1912'
$ws1.Cells.Item(354, 1).Value = '15:45:31'
$ws1.Cells.Item(354, 2).Value = '17:35'
$ws1.Cells.Item(354, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(354, 4).Value = 110
$ws1.Cells.Item(354, 5).Value = 'LP1912'
$ws1.Cells.Item(355, 1).Value = '16:13:19'
$ws1.Cells.Item(355, 2).Value = '17:35'
$ws1.Cells.Item(355, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(355, 4).Value = 82
$ws1.Cells.Item(355, 5).Value = 'LP1912'
$ws1.Cells.Item(356, 1).Value = '15:45:31'
$ws1.Cells.Item(356, 2).Value = '17:36'
$ws1.Cells.Item(356, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(356, 4).Value = 111
$ws1.Cells.Item(356, 5).Value = 'LP1912'
$ws1.Cells.Item(357, 1).Value = '15:45:31'
$ws1.Cells.Item(357, 2).Value = '17:38'
$ws1.Cells.Item(357, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(357, 4).Value = 113
$ws1.Cells.Item(357, 5).Value = 'LP1912'
$ws1.Cells.Item(358, 1).Value = '16:27:37'
$ws1.Cells.Item(358, 2).Value = '17:44'
$ws1.Cells.Item(358, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(358, 4).Value = 77
$ws1.Cells.Item(358, 5).Value = 'LP1912'
$ws1.Cells.Item(359, 1).Value = '15:57:19'
$ws1.Cells.Item(359, 2).Value = '17:45'
$ws1.Cells.Item(359, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(359, 4).Value = 108
$ws1.Cells.Item(359, 5).Value = 'LP1912'
$ws1.Cells.Item(360, 1).Value = '16:27:37'
$ws1.Cells.Item(360, 2).Value = '17:49'
$ws1.Cells.Item(360, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(360, 4).Value = 82
$ws1.Cells.Item(360, 5).Value = 'LP1912'
$ws1.Cells.Item(361, 1).Value = '15:57:19'
$ws1.Cells.Item(361, 2).Value = '17:51'
$ws1.Cells.Item(361, 3).Value = '215_EL PELIGRO'
$ws1.Cells.Item(361, 4).Value = 114
$ws1.Cells.Item(361, 5).Value = 'LP1912'
$ws1.Cells.Item(362, 1).Value = '16:27:37'
$ws1.Cells.Item(362, 2).Value = '18:02'
$ws1.Cells.Item(362, 3).Value = '17_ROMERO'
$ws1.Cells.Item(362, 4).Value = 95
$ws1.Cells.Item(362, 5).Value = 'LP1912'
$ws1.Cells.Item(363, 1).Value = '16:13:19'
$ws1.Cells.Item(363, 2).Value = '18:03'
$ws1.Cells.Item(363, 3).Value = '17_ROMERO'
$ws1.Cells.Item(363, 4).Value = 110
$ws1.Cells.Item(363, 5).Value = 'LP1912'
$ws1.Cells.Item(364, 1).Value = '16:36:34'
$ws1.Cells.Item(364, 2).Value = '18:04'
$ws1.Cells.Item(364, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(364, 4).Value = 88
$ws1.Cells.Item(364, 5).Value = 'LP1912'
$ws1.Cells.Item(365, 1).Value = '16:13:19'
$ws1.Cells.Item(365, 2).Value = '18:04'
$ws1.Cells.Item(365, 3).Value = '14_ABASTO'
$ws1.Cells.Item(365, 4).Value = 111
$ws1.Cells.Item(365, 5).Value = 'LP1912'
$ws1.Cells.Item(366, 1).Value = '16:27:37'
$ws1.Cells.Item(366, 2).Value = '18:24'
$ws1.Cells.Item(366, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(366, 4).Value = 117
$ws1.Cells.Item(366, 5).Value = 'LP1912'
$ws1.Cells.Item(367, 1).Value = '16:36:34'
$ws1.Cells.Item(367, 2).Value = '18:34'
$ws1.Cells.Item(367, 3).Value = '14X44_ABASTO'
$ws1.Cells.Item(367, 4).Value = 118
$ws1.Cells.Item(367, 5).Value = 'LP1912'
